$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all updated cells keep their original text formatting (avoid Excel
# auto-converting numeric-looking strings like "1.00" or "0.0240" into real
# numbers). Each cell is formatted individually because this COM runtime only
# applies NumberFormat to the first area of a multi-area (comma) range.
$textCells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7",
    "E7", "D8", "E8", "E9", "E10", "D11", "E11", "E12", "D13", "E13",
    "D14", "E14", "D15", "E15", "E16", "D17", "E17", "D18", "E18", "D19",
    "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24",
    "E24", "E25", "E26", "D27", "E27", "D28", "E28", "B29", "C29", "D29",
    "E29", "B30", "C30", "D30", "E30", "E31", "B32", "C32", "D32", "E32",
    "B33", "C33", "D33", "E33", "D34", "E34", "D35", "E35", "E36", "D37",
    "E37", "D38", "E38", "D39", "E39", "D40", "E40", "B41", "C41", "D41",
    "E41", "B42", "C42", "D42", "E42", "D43", "E43", "B44", "C44", "D44",
    "E44", "E45", "B46", "C46", "D46", "E46", "B47", "C47", "D47", "E47",
    "B48", "C48", "D48", "E48", "B49", "C49", "D49", "E49", "D50", "E50",
    "D51", "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated coin price / volume(1h) figures, and the rows whose coin ranking
# order swapped (Aptos/PancakeSwap, Monero/USDe, OKB/Filecoin,
# Maker/InjectiveProtocol/WhiteBITCoin/Hedera/EnergySwap).

# Row 2
$ws.Range("D2").Value = '61.974.68'
$ws.Range("E2").Value = '  -2.29%  '

# Row 3
$ws.Range("D3").Value = '2.587.42'
$ws.Range("E3").Value = '  -4.65%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = '553.18'
$ws.Range("E5").Value = '  -1.40%  '

# Row 6
$ws.Range("D6").Value = '154.97'
$ws.Range("E6").Value = '  -1.63%  '

# Row 7
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.11%  '

# Row 8
$ws.Range("D8").Value = '0.594'
$ws.Range("E8").Value = '  +0.65%  '

# Row 9
$ws.Range("E9").Value = '  -2.74%  '

# Row 10
$ws.Range("E10").Value = '  -3.53%  '

# Row 11
$ws.Range("D11").Value = '5.43'
$ws.Range("E11").Value = '  -2.24%  '

# Row 12
$ws.Range("E12").Value = '  -2.11%  '

# Row 13
$ws.Range("D13").Value = '3.040.57'
$ws.Range("E13").Value = '  -4.80%  '

# Row 14
$ws.Range("D14").Value = '25.62'
$ws.Range("E14").Value = '  -3.25%  '

# Row 15
$ws.Range("D15").Value = '61.822.05'
$ws.Range("E15").Value = '  -2.31%  '

# Row 16
$ws.Range("E16").Value = '  -2.77%  '

# Row 17
$ws.Range("D17").Value = '2.589.24'
$ws.Range("E17").Value = '  -4.66%  '

# Row 18
$ws.Range("D18").Value = '11.62'
$ws.Range("E18").Value = '  -4.42%  '

# Row 19
$ws.Range("D19").Value = '4.54'
$ws.Range("E19").Value = '  -2.81%  '

# Row 20
$ws.Range("D20").Value = '339.12'
$ws.Range("E20").Value = '  -3.48%  '

# Row 21
$ws.Range("D21").Value = '6.04'
$ws.Range("E21").Value = '  -6.33%  '

# Row 22
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  +0.04%  '

# Row 23
$ws.Range("D23").Value = '0.498'
$ws.Range("E23").Value = '  -2.64%  '

# Row 24
$ws.Range("D24").Value = '62.49'
$ws.Range("E24").Value = '  -2.79%  '

# Row 25
$ws.Range("E25").Value = '  -0.32%  '

# Row 26
$ws.Range("E26").Value = '  -0.04%  '

# Row 27
$ws.Range("D27").Value = '8.02'
$ws.Range("E27").Value = '  -2.35%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0836'
$ws.Range("E28").Value = '  -5.94%  '

# Row 29
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '1.91'
$ws.Range("E29").Value = '  -2.61%  '

# Row 30
$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").Value = '7.04'
$ws.Range("E30").Value = '  -1.83%  '

# Row 31
$ws.Range("E31").Value = '  -4.32%  '

# Row 32
$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.03%  '

# Row 33
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").Value = '159.37'
$ws.Range("E33").Value = '  -4.29%  '

# Row 34
$ws.Range("D34").Value = '19.25'
$ws.Range("E34").Value = '  -2.75%  '

# Row 35
$ws.Range("D35").Value = '4.68'
$ws.Range("E35").Value = '  -3.32%  '

# Row 36
$ws.Range("E36").Value = '  -4.77%  '

# Row 37
$ws.Range("D37").Value = '1.75'
$ws.Range("E37").Value = '  -1.42%  '

# Row 38
$ws.Range("D38").Value = '337.88'
$ws.Range("E38").Value = '  -2.61%  '

# Row 39
$ws.Range("D39").Value = '6.01'
$ws.Range("E39").Value = '  -2.21%  '

# Row 40
$ws.Range("D40").Value = '0.895'
$ws.Range("E40").Value = '  -6.79%  '

# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '3.90'
$ws.Range("E41").Value = '  -3.28%  '

# Row 42
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '37.50'
$ws.Range("E42").Value = '  -2.39%  '

# Row 43
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  +0.08%  '

# Row 44
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '20.41'
$ws.Range("E44").Value = '  -4.67%  '

# Row 45
$ws.Range("E45").Value = '  -2.63%  '

# Row 46
$ws.Range("B46").Value = 'WhiteBITCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D46").Value = '10.94'
$ws.Range("E46").Value = '  -1.14%  '

# Row 47
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.114.95'
$ws.Range("E47").Value = '  +0.74%  '

# Row 48
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '0.0547'
$ws.Range("E48").Value = '  -4.55%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '19.63'
$ws.Range("E49").Value = '  -5.03%  '

# Row 50
$ws.Range("D50").Value = '0.0964'
$ws.Range("E50").Value = '  -2.18%  '

# Row 51
$ws.Range("D51").Value = '0.0240'
$ws.Range("E51").Value = '  -2.64%  '
